# edit.ps1 - reproduces the authored change:
#   * slide 1: the "A" / "B" caption textboxes become lower-case "a" / "b"
#     and their (auto-fit) boxes shrink to the narrower glyph width
#   * the cached "datetimeFigureOut" footer field (slide master + every
#     slide layout) is re-stamped from 2/28/20 to 3/2/20, exactly like
#     PowerPoint does whenever the deck is re-saved on a later date

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 caption textboxes: "A" -> "a" and "B" -> "b", with the
#    spAutoFit textbox shrinking to the new (narrower) glyph width.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$boxA = $slide.Shapes.Item(3)   # "TextBox 7" - currently "A"
$boxA.TextFrame.TextRange.Text = "a"
$boxA.Width = (287258 + 0.5) / 12700

$boxB = $slide.Shapes.Item(4)   # "TextBox 8" - currently "B"
$boxB.TextFrame.TextRange.Text = "b"
$boxB.Width = (300082 + 0.5) / 12700

# ---------------------------------------------------------------------
# 2) Re-cache the datetimeFigureOut date placeholder text: 2/28/20 -> 3/2/20
#    on the slide master and on every slide layout.
# ---------------------------------------------------------------------
$newDate = "3/2/20"

$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    $layoutShapes = $layout.Shapes
    for ($i = 1; $i -le $layoutShapes.Count; $i++) {
        $sh = $layoutShapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
